$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update cell E8 from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Set selection to E8
$ws.Range("E8").Select()
